# Auto-generated COM-interop script implementing the 2022-Q4 sheet edit.
$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 25
$summary.Cells.Item(2, 4).Value = 5.35
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 12
$summary.Cells.Item(3, 4).Value = 3.95
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q2"
$summary.Cells.Item(4, 3).Value = 3
$summary.Cells.Item(4, 4).Value = 0.8
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2022-Q1"
$summary.Cells.Item(5, 3).Value = 9
$summary.Cells.Item(5, 4).Value = 1.38
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q4"
$summary.Cells.Item(6, 3).Value = 8
$summary.Cells.Item(6, 4).Value = 1.62
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(7, 2).Value = "2021-Q3"
$summary.Cells.Item(7, 3).Value = 11
$summary.Cells.Item(7, 4).Value = 0.38
$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(8, 2).Value = "2021-Q2"
$summary.Cells.Item(8, 3).Value = 3
$summary.Cells.Item(8, 4).Value = 1.17
$summary.Cells.Item(9, 1).Value = 7
$summary.Cells.Item(9, 2).Value = "2021-Q1"
$summary.Cells.Item(9, 3).Value = 4
$summary.Cells.Item(9, 4).Value = 1.09

$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)
$summary.Range("A1").Select()

$before = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($before)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Copy the header style (bold/border, style index 2) from the
# neighbouring "2022-Q3" sheet's identical header row.
$before.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Data rows: column A (index) and column H (rank) are numeric;
# columns B-G are textual (fund code/name/percentages kept as literal
# text, matching the source sheet's inlineStr cells) -- written via a
# text-literal formula then flattened to static values with PasteSpecial
# so no formula/style residue is left behind.
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Formula = "=`"166301`""
$q4.Cells.Item(2, 3).Formula = "=`"华商新趋势优选灵活配置混合`""
$q4.Cells.Item(2, 4).Formula = "=`"98.72`""
$q4.Cells.Item(2, 5).Formula = "=`"74.44`""
$q4.Cells.Item(2, 6).Formula = "=`"1.91`""
$q4.Cells.Item(2, 7).Formula = "=`"1.8856`""
$q4.Cells.Item(2, 8).Value = 5
$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Formula = "=`"000390`""
$q4.Cells.Item(3, 3).Formula = "=`"华商优势行业混合`""
$q4.Cells.Item(3, 4).Formula = "=`"35.74`""
$q4.Cells.Item(3, 5).Formula = "=`"90.63`""
$q4.Cells.Item(3, 6).Formula = "=`"2.13`""
$q4.Cells.Item(3, 7).Formula = "=`"0.7613`""
$q4.Cells.Item(3, 8).Value = 8
$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Formula = "=`"010761`""
$q4.Cells.Item(4, 3).Formula = "=`"华商甄选回报混合A`""
$q4.Cells.Item(4, 4).Formula = "=`"40.49`""
$q4.Cells.Item(4, 5).Formula = "=`"70.56`""
$q4.Cells.Item(4, 6).Formula = "=`"1.38`""
$q4.Cells.Item(4, 7).Formula = "=`"0.5588`""
$q4.Cells.Item(4, 8).Value = 7
$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Formula = "=`"630002`""
$q4.Cells.Item(5, 3).Formula = "=`"华商盛世成长混合`""
$q4.Cells.Item(5, 4).Formula = "=`"23.51`""
$q4.Cells.Item(5, 5).Formula = "=`"90.49`""
$q4.Cells.Item(5, 6).Formula = "=`"2.19`""
$q4.Cells.Item(5, 7).Formula = "=`"0.5149`""
$q4.Cells.Item(5, 8).Value = 10
$q4.Cells.Item(6, 1).Value = 4
$q4.Cells.Item(6, 2).Formula = "=`"004046`""
$q4.Cells.Item(6, 3).Formula = "=`"华夏新锦顺灵活配置混合A`""
$q4.Cells.Item(6, 4).Formula = "=`"24.34`""
$q4.Cells.Item(6, 5).Formula = "=`"71.43`""
$q4.Cells.Item(6, 6).Formula = "=`"1.94`""
$q4.Cells.Item(6, 7).Formula = "=`"0.4722`""
$q4.Cells.Item(6, 8).Value = 2
$q4.Cells.Item(7, 1).Value = 5
$q4.Cells.Item(7, 2).Formula = "=`"013627`""
$q4.Cells.Item(7, 3).Formula = "=`"华夏周期驱动混合C`""
$q4.Cells.Item(7, 4).Formula = "=`"4.87`""
$q4.Cells.Item(7, 5).Formula = "=`"86.45`""
$q4.Cells.Item(7, 6).Formula = "=`"5.97`""
$q4.Cells.Item(7, 7).Formula = "=`"0.2907`""
$q4.Cells.Item(7, 8).Value = 2
$q4.Cells.Item(8, 1).Value = 6
$q4.Cells.Item(8, 2).Formula = "=`"016049`""
$q4.Cells.Item(8, 3).Formula = "=`"华商甄选回报混合C`""
$q4.Cells.Item(8, 4).Formula = "=`"20.87`""
$q4.Cells.Item(8, 5).Formula = "=`"70.56`""
$q4.Cells.Item(8, 6).Formula = "=`"1.38`""
$q4.Cells.Item(8, 7).Formula = "=`"0.2880`""
$q4.Cells.Item(8, 8).Value = 7
$q4.Cells.Item(9, 1).Value = 7
$q4.Cells.Item(9, 2).Formula = "=`"013626`""
$q4.Cells.Item(9, 3).Formula = "=`"华夏周期驱动混合A`""
$q4.Cells.Item(9, 4).Formula = "=`"4.12`""
$q4.Cells.Item(9, 5).Formula = "=`"86.45`""
$q4.Cells.Item(9, 6).Formula = "=`"5.97`""
$q4.Cells.Item(9, 7).Formula = "=`"0.2460`""
$q4.Cells.Item(9, 8).Value = 2
$q4.Cells.Item(10, 1).Value = 8
$q4.Cells.Item(10, 2).Formula = "=`"005660`""
$q4.Cells.Item(10, 3).Formula = "=`"嘉实资源精选股票A`""
$q4.Cells.Item(10, 4).Formula = "=`"2.63`""
$q4.Cells.Item(10, 5).Formula = "=`"93.36`""
$q4.Cells.Item(10, 6).Formula = "=`"4.43`""
$q4.Cells.Item(10, 7).Formula = "=`"0.1165`""
$q4.Cells.Item(10, 8).Value = 10
$q4.Cells.Item(11, 1).Value = 9
$q4.Cells.Item(11, 2).Formula = "=`"008488`""
$q4.Cells.Item(11, 3).Formula = "=`"华商恒益稳健混合`""
$q4.Cells.Item(11, 4).Formula = "=`"4.22`""
$q4.Cells.Item(11, 5).Formula = "=`"49.85`""
$q4.Cells.Item(11, 6).Formula = "=`"1.52`""
$q4.Cells.Item(11, 7).Formula = "=`"0.0641`""
$q4.Cells.Item(11, 8).Value = 7
$q4.Cells.Item(12, 1).Value = 10
$q4.Cells.Item(12, 2).Formula = "=`"005661`""
$q4.Cells.Item(12, 3).Formula = "=`"嘉实资源精选股票C`""
$q4.Cells.Item(12, 4).Formula = "=`"1.06`""
$q4.Cells.Item(12, 5).Formula = "=`"93.36`""
$q4.Cells.Item(12, 6).Formula = "=`"4.43`""
$q4.Cells.Item(12, 7).Formula = "=`"0.0470`""
$q4.Cells.Item(12, 8).Value = 10
$q4.Cells.Item(13, 1).Value = 11
$q4.Cells.Item(13, 2).Formula = "=`"004047`""
$q4.Cells.Item(13, 3).Formula = "=`"华夏新锦顺灵活配置混合C`""
$q4.Cells.Item(13, 4).Formula = "=`"1.24`""
$q4.Cells.Item(13, 5).Formula = "=`"71.43`""
$q4.Cells.Item(13, 6).Formula = "=`"1.94`""
$q4.Cells.Item(13, 7).Formula = "=`"0.0241`""
$q4.Cells.Item(13, 8).Value = 2
$q4.Cells.Item(14, 1).Value = 12
$q4.Cells.Item(14, 2).Formula = "=`"005161`""
$q4.Cells.Item(14, 3).Formula = "=`"华商上游产业股票`""
$q4.Cells.Item(14, 4).Formula = "=`"0.55`""
$q4.Cells.Item(14, 5).Formula = "=`"88.76`""
$q4.Cells.Item(14, 6).Formula = "=`"4.22`""
$q4.Cells.Item(14, 7).Formula = "=`"0.0232`""
$q4.Cells.Item(14, 8).Value = 4
$q4.Cells.Item(15, 1).Value = 13
$q4.Cells.Item(15, 2).Formula = "=`"006401`""
$q4.Cells.Item(15, 3).Formula = "=`"先锋量化优选灵活配置混合A`""
$q4.Cells.Item(15, 4).Formula = "=`"0.59`""
$q4.Cells.Item(15, 5).Formula = "=`"92.91`""
$q4.Cells.Item(15, 6).Formula = "=`"2.63`""
$q4.Cells.Item(15, 7).Formula = "=`"0.0155`""
$q4.Cells.Item(15, 8).Value = 1
$q4.Cells.Item(16, 1).Value = 14
$q4.Cells.Item(16, 2).Formula = "=`"011888`""
$q4.Cells.Item(16, 3).Formula = "=`"民生加银周期优选混合型证券投资基金A`""
$q4.Cells.Item(16, 4).Formula = "=`"0.35`""
$q4.Cells.Item(16, 5).Formula = "=`"92.72`""
$q4.Cells.Item(16, 6).Formula = "=`"4.01`""
$q4.Cells.Item(16, 7).Formula = "=`"0.0140`""
$q4.Cells.Item(16, 8).Value = 3
$q4.Cells.Item(17, 1).Value = 15
$q4.Cells.Item(17, 2).Formula = "=`"008629`""
$q4.Cells.Item(17, 3).Formula = "=`"大成景瑞稳健配置混合A`""
$q4.Cells.Item(17, 4).Formula = "=`"0.66`""
$q4.Cells.Item(17, 5).Formula = "=`"29.66`""
$q4.Cells.Item(17, 6).Formula = "=`"1.42`""
$q4.Cells.Item(17, 7).Formula = "=`"0.0094`""
$q4.Cells.Item(17, 8).Value = 9
$q4.Cells.Item(18, 1).Value = 16
$q4.Cells.Item(18, 2).Formula = "=`"519172`""
$q4.Cells.Item(18, 3).Formula = "=`"浦银安盛睿智精选灵活配置混合A`""
$q4.Cells.Item(18, 4).Formula = "=`"0.21`""
$q4.Cells.Item(18, 5).Formula = "=`"89.04`""
$q4.Cells.Item(18, 6).Formula = "=`"2.66`""
$q4.Cells.Item(18, 7).Formula = "=`"0.0056`""
$q4.Cells.Item(18, 8).Value = 10
$q4.Cells.Item(19, 1).Value = 17
$q4.Cells.Item(19, 2).Formula = "=`"006402`""
$q4.Cells.Item(19, 3).Formula = "=`"先锋量化优选灵活配置混合C`""
$q4.Cells.Item(19, 4).Formula = "=`"0.21`""
$q4.Cells.Item(19, 5).Formula = "=`"92.91`""
$q4.Cells.Item(19, 6).Formula = "=`"2.63`""
$q4.Cells.Item(19, 7).Formula = "=`"0.0055`""
$q4.Cells.Item(19, 8).Value = 1
$q4.Cells.Item(20, 1).Value = 18
$q4.Cells.Item(20, 2).Formula = "=`"008630`""
$q4.Cells.Item(20, 3).Formula = "=`"大成景瑞稳健配置混合C`""
$q4.Cells.Item(20, 4).Formula = "=`"0.25`""
$q4.Cells.Item(20, 5).Formula = "=`"29.66`""
$q4.Cells.Item(20, 6).Formula = "=`"1.42`""
$q4.Cells.Item(20, 7).Formula = "=`"0.0036`""
$q4.Cells.Item(20, 8).Value = 9
$q4.Cells.Item(21, 1).Value = 19
$q4.Cells.Item(21, 2).Formula = "=`"011889`""
$q4.Cells.Item(21, 3).Formula = "=`"民生加银周期优选混合型证券投资基金C`""
$q4.Cells.Item(21, 4).Formula = "=`"0.09`""
$q4.Cells.Item(21, 5).Formula = "=`"92.72`""
$q4.Cells.Item(21, 6).Formula = "=`"4.01`""
$q4.Cells.Item(21, 7).Formula = "=`"0.0036`""
$q4.Cells.Item(21, 8).Value = 3
$q4.Cells.Item(22, 1).Value = 20
$q4.Cells.Item(22, 2).Formula = "=`"519173`""
$q4.Cells.Item(22, 3).Formula = "=`"浦银安盛睿智精选灵活配置混合C`""
$q4.Cells.Item(22, 4).Formula = "=`"0.13`""
$q4.Cells.Item(22, 5).Formula = "=`"89.04`""
$q4.Cells.Item(22, 6).Formula = "=`"2.66`""
$q4.Cells.Item(22, 7).Formula = "=`"0.0035`""
$q4.Cells.Item(22, 8).Value = 10
$q4.Cells.Item(23, 1).Value = 21
$q4.Cells.Item(23, 2).Formula = "=`"004727`""
$q4.Cells.Item(23, 3).Formula = "=`"先锋聚优灵活配置混合C`""
$q4.Cells.Item(23, 4).Formula = "=`"0.02`""
$q4.Cells.Item(23, 5).Formula = "=`"93.09`""
$q4.Cells.Item(23, 6).Formula = "=`"2.55`""
$q4.Cells.Item(23, 7).Formula = "=`"0.0005`""
$q4.Cells.Item(23, 8).Value = 5
$q4.Cells.Item(24, 1).Value = 22
$q4.Cells.Item(24, 2).Formula = "=`"003587`""
$q4.Cells.Item(24, 3).Formula = "=`"先锋精一灵活配置混合C`""
$q4.Cells.Item(24, 4).Formula = "=`"0.02`""
$q4.Cells.Item(24, 5).Formula = "=`"93.29`""
$q4.Cells.Item(24, 6).Formula = "=`"2.66`""
$q4.Cells.Item(24, 7).Formula = "=`"0.0005`""
$q4.Cells.Item(24, 8).Value = 3
$q4.Cells.Item(25, 1).Value = 23
$q4.Cells.Item(25, 2).Formula = "=`"004726`""
$q4.Cells.Item(25, 3).Formula = "=`"先锋聚优灵活配置混合A`""
$q4.Cells.Item(25, 4).Formula = "=`"0.01`""
$q4.Cells.Item(25, 5).Formula = "=`"93.09`""
$q4.Cells.Item(25, 6).Formula = "=`"2.55`""
$q4.Cells.Item(25, 7).Formula = "=`"0.0003`""
$q4.Cells.Item(25, 8).Value = 5
$q4.Cells.Item(26, 1).Value = 24
$q4.Cells.Item(26, 2).Formula = "=`"003586`""
$q4.Cells.Item(26, 3).Formula = "=`"先锋精一灵活配置混合A`""
$q4.Cells.Item(26, 4).Formula = "=`"0.01`""
$q4.Cells.Item(26, 5).Formula = "=`"93.29`""
$q4.Cells.Item(26, 6).Formula = "=`"2.66`""
$q4.Cells.Item(26, 7).Formula = "=`"0.0003`""
$q4.Cells.Item(26, 8).Value = 3

$textRange = $q4.Range("B2:G26")
$textRange.Copy()
$textRange.PasteSpecial(-4163)

$summary.Range("A2").Copy()
$aRange = $q4.Range("A2:A26")
$aRange.PasteSpecial(-4122)

$q4.Range("A1").Select()

